$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3044
$ws1.Range("F3").Value = 469
$ws1.Range("F7").Value = 219
$ws1.Range("F8").Value = 14573
$ws1.Range("F9").Value = 164
$ws1.Range("F10").Value = 124
$ws1.Range("F11").Value = 5818
$ws1.Range("F12").Value = 592
$ws1.Range("F13").Value = 74
$ws1.Range("F19").Value = 185
$ws1.Range("F20").Value = 798
$ws1.Range("F21").Value = 2941
$ws1.Range("F22").Value = 57
$ws1.Range("F23").Value = 10606
$ws1.Range("F26").Value = 85

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3044
$ws4.Range("F4").Value = 469
$ws4.Range("F8").Value = 219
$ws4.Range("F9").Value = 14573
$ws4.Range("F10").Value = 164
$ws4.Range("F11").Value = 124
$ws4.Range("F12").Value = 5818
$ws4.Range("F13").Value = 592
$ws4.Range("F14").Value = 74
$ws4.Range("F20").Value = 185
$ws4.Range("F21").Value = 798
$ws4.Range("F22").Value = 2941
$ws4.Range("F23").Value = 57
$ws4.Range("F25").Value = 10606
$ws4.Range("F28").Value = 85
